$newRows = @(
    ,@('n8n×Python×AIで公開レポート自動探索・抽出・分類・登録フロー(PoC開発)','システム開発','100,000 円 ~ 200,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5413825',538,'🔥AI,Python ◆開発')
    ,@('【急募】AI×LINE開発をリード!医療機関向けアプリのサーバーサイドエンジニア募集(フルリモート)','システム開発','500,000 円 ~ 1,000,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5413230',385,'🔥AI,Ai ◆開発 ◇アプリ')
    ,@('【急募】不動産評価システムのAIチャットポット開発者募集','システム開発','100,000 円 ~ 200,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5413280',383,'🔥AI,Ai ◆開発')
    ,@('【急募】ローカルAI開発プロジェクトの協力者を探しています!','システム開発','300,000 円 ~ 500,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5413402',375,'🔥AI,Ai ◆開発')
    ,@('【高報酬/リモート可/法人可】グローバルHRベンチャーでAIを活用し業務効率化を推進してくださる方!','システム開発','5,000 円 ~ 10,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5413210',370,'🔥AI,Ai ◆効率化')
    ,@('【急募】AIテキスト抜粋アプリのプロンプト最適化依頼','システム開発','100,000 円 ~ 200,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5413215',338,'🔥AI,Ai ◇アプリ')
    ,@('【Azure/RAG】社内文書検索AIチャットボットの精度向上&内製化支援パートナー募集!','システム開発','300,000 円 ~ 500,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5413954',310,'🔥AI,Ai')
    ,@('生成AIの技術顧問を募集!事業の技術選定をリードするAI専門家を募集! 【週1日〜/フルリモート】','システム開発','200,000 円 ~ 300,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5413955',303,'🔥AI,Ai')
    ,@('コスパスポーツジムの高速自動予約botの開発','システム開発','50,000 円 ~ 100,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5413835',173,'★bot ◆開発')
    ,@('3Dプリント用データのWeb自動チェック&変換&カラー補正ツール|開発パートナー募集','システム開発','5,000,000 円 ~ / 固定','期限情報なし','https://www.lancers.jp/work/detail/5413508',135,'◆ツール,開発')
    ,@('【長期依頼】海外クリエイター向けサービスの開発保守|Laravel+Livewireエンジニア募集','システム開発','300,000 円 ~ 500,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5414105',75,'◆開発')
    ,@('【急募】見積書自動作成機能の開発','システム開発','500,000 円 ~ 1,000,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5414108',75,'◆開発')
    ,@('セレニウムを用いた自動発注ツールの修正・機能追加','システム開発','50,000 円 ~ 100,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5413916',73,'◆ツール')
    ,@('IB報酬を得るための高性能EA開発依頼','システム開発','100,000 円 ~ 200,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5413293',68,'◆開発')
    ,@('【音声コマンド起動】超小型・低電力レコーダーのプロトタイプ開発','システム開発','20,000 円 ~ 50,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5413958',63,'◆開発')
    ,@('wordpressレンダリングを妨げるリソースの除外','システム開発','200,000 円 ~ 300,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5016989',33,'○WordPress')
    ,@('【急募】Cloud RunでWordPress構築のプロを探しています!','システム開発','10,000 円 ~ 20,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5413043',25,'○WordPress')
    ,@('【急募】16タイプ診断コンテンツのLP制作','システム開発','300,000 円 ~ 500,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5408735',25,'')
    ,@('【急募】Teams Roomsの設定設置と保守サポート依頼','システム開発','100,000 円 ~ 200,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5408814',18,'')
    ,@('Access 32bitから64bitへの修正改善依頼','システム開発','5,000 円 ~ 10,000 円 / 固定','期限情報なし','https://www.lancers.jp/work/detail/5413333',10,'')
)

$timestamp = "2025-10-16 01:45:48"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove existing hyperlinks; they will be rebuilt after the row insert so
# that every F-column cell ends up pointing at the right URL again.
$ws.Hyperlinks.Delete()

# Insert two new rows right after row 11 (before the old row 12), which
# shifts the old rows 12-19 down to 14-21.
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $i + 2
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $data[0]
    $ws.Cells.Item($r, 3).Value = $data[1]
    $ws.Cells.Item($r, 4).Value = $data[2]
    $ws.Cells.Item($r, 5).Value = $data[3]
    $ws.Cells.Item($r, 6).Value = $data[4]
    $ws.Cells.Item($r, 7).Value = $data[5]
    if ($data[6] -ne '') {
        $ws.Cells.Item($r, 8).Value = $data[6]
    } else {
        $ws.Cells.Item($r, 8).Value = ""
    }

    $ws.Hyperlinks.Add($ws.Range("F" + $r), $data[4])
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
}

$ws.Range("A1").Select() | Out-Null
